$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Raw data table (A:D) -------------------------------------------------
# Existing rows are 1 (header) .. 11. Append three more measurements.
$newRows = @(
    @(800,  1.98, 0.84, -23),
    @(1000, 1.98, 0.82, -14),
    @(1200, 1.98, 0.82, -10)
)

$r = 12
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# ---- Plotting helper table (H:J) ------------------------------------------
# H = freq (copy of A), I = Vout/Vin (= C/B), J = Phase (copy of D)
$allRows = @(
    @(5,    1.98, 6.32, -9),
    @(8,    1.98, 6.16, -11),
    @(14,   1.96, 5.76, -20),
    @(23,   1.94, 5.12, -32),
    @(39,   1.96, 4.16, -39),
    @(65,   1.96, 2.96, -51),
    @(108,  1.96, 2.04, -50),
    @(180,  1.98, 1.42, -51),
    @(300,  1.98, 1.08, -39),
    @(500,  1.96, 0.90, -25),
    @(800,  1.98, 0.84, -23),
    @(1000, 1.98, 0.82, -14),
    @(1200, 1.98, 0.82, -10)
)

$r = 2
foreach ($row in $allRows) {
    $ws.Cells.Item($r, 8).Value = $row[0]
    $ws.Cells.Item($r, 10).Value = $row[3]
    $r = $r + 1
}

$ws.Cells.Item(2, 9).Formula = "=C2/B2"
$ws.Range("I3:I14").Formula = "=C3/B3"

# ---- Sheet view / selection -------------------------------------------------
$ws.Range("S15").Select()

# ---- Page setup --------------------------------------------------------
$ws.PageSetup.Orientation = 1
